# Daily attendance processing - 2025-10-11 10:19:20
#
# Normalizes the "Recorded By" attendee list ordering (column G) on the
# "Session Analysis Results" sheet so the "System" entry is listed last
# instead of first (e.g. "System, dnasr281@gmail.com" ->
# "dnasr281@gmail.com, System"), and refreshes a handful of attendance
# counts / percentages (columns H and S) to reflect corrected headcounts.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$ws.Range("G2").Value = "backup@backdoor.com, system, System"
$ws.Range("G3").Value = "dnasr281@gmail.com, System"
$ws.Range("G4").Value = "backup@backdoor.com, System"
$ws.Range("G5").Value = "backup@backdoor.com, System"
$ws.Range("G6").Value = "dnasr281@gmail.com, System"
$ws.Range("G7").Value = "admin@admin.com, System"
$ws.Range("G10").Value = "dnasr281@gmail.com, System"
$ws.Range("G11").Value = "dnasr281@gmail.com, System"
$ws.Range("G12").Value = "dnasr281@gmail.com, System"
$ws.Range("G13").Value = "dnasr281@gmail.com, System"
$ws.Range("G14").Value = "dnasr281@gmail.com, System"
$ws.Range("G15").Value = "dnasr281@gmail.com, System"
$ws.Range("S15").Value = "'69.8%"
$ws.Range("S16").Value = "'66.0%"
$ws.Range("G17").Value = "dnasr281@gmail.com, System"
$ws.Range("H17").Value = "38/53"
$ws.Range("S17").Value = "'61.6%"
$ws.Range("G29").Value = "backup@backdoor.com, system, System"
$ws.Range("G30").Value = "dnasr281@gmail.com, System"
$ws.Range("G32").Value = "backup@backdoor.com, System"
$ws.Range("G33").Value = "dnasr281@gmail.com, System"
$ws.Range("G34").Value = "admin@admin.com, System"
$ws.Range("G37").Value = "dnasr281@gmail.com, System"
$ws.Range("G38").Value = "dnasr281@gmail.com, System"
$ws.Range("G39").Value = "dnasr281@gmail.com, System"
$ws.Range("G40").Value = "dnasr281@gmail.com, System"
$ws.Range("G41").Value = "dnasr281@gmail.com, System"
$ws.Range("G42").Value = "dnasr281@gmail.com, System"
$ws.Range("G44").Value = "dnasr281@gmail.com, System"
$ws.Range("H44").Value = "43/56"
$ws.Range("G56").Value = "backup@backdoor.com, system, System"
$ws.Range("G57").Value = "dnasr281@gmail.com, System"
$ws.Range("H57").Value = "28/55"
$ws.Range("G58").Value = "backup@backdoor.com, System"
$ws.Range("G59").Value = "backup@backdoor.com, System"
$ws.Range("G60").Value = "dnasr281@gmail.com, System"
$ws.Range("H60").Value = "37/55"
$ws.Range("G61").Value = "admin@admin.com, System"
$ws.Range("G64").Value = "dnasr281@gmail.com, System"
$ws.Range("G65").Value = "dnasr281@gmail.com, System"
$ws.Range("G66").Value = "dnasr281@gmail.com, System"
$ws.Range("G67").Value = "dnasr281@gmail.com, System"
$ws.Range("G68").Value = "dnasr281@gmail.com, System"
$ws.Range("G69").Value = "dnasr281@gmail.com, System"
$ws.Range("G71").Value = "dnasr281@gmail.com, System"
$ws.Range("G84").Value = "backup@backdoor.com, System"
$ws.Range("G85").Value = "backup@backdoor.com, System"
$ws.Range("G86").Value = "dnasr281@gmail.com, System"
$ws.Range("G87").Value = "dnasr281@gmail.com, System"
$ws.Range("G88").Value = "dnasr281@gmail.com, System"
$ws.Range("G89").Value = "dnasr281@gmail.com, System"
$ws.Range("G90").Value = "dnasr281@gmail.com, admin@admin.com"
$ws.Range("G93").Value = "dnasr281@gmail.com, System"
$ws.Range("G95").Value = "dnasr281@gmail.com, System"
$ws.Range("G96").Value = "dnasr281@gmail.com, System"
$ws.Range("G110").Value = "backup@backdoor.com, System"
$ws.Range("G111").Value = "backup@backdoor.com, System"
$ws.Range("G112").Value = "dnasr281@gmail.com, System"
$ws.Range("G113").Value = "dnasr281@gmail.com, System"
$ws.Range("G114").Value = "dnasr281@gmail.com, System"
$ws.Range("G115").Value = "dnasr281@gmail.com, System"
$ws.Range("G116").Value = "dnasr281@gmail.com, admin@admin.com"
$ws.Range("G119").Value = "dnasr281@gmail.com, System"
$ws.Range("G121").Value = "dnasr281@gmail.com, System"
$ws.Range("G122").Value = "dnasr281@gmail.com, System"
$ws.Range("G136").Value = "backup@backdoor.com, System"
$ws.Range("G137").Value = "backup@backdoor.com, System"
$ws.Range("G138").Value = "dnasr281@gmail.com, System"
$ws.Range("G139").Value = "dnasr281@gmail.com, System"
$ws.Range("G140").Value = "dnasr281@gmail.com, System"
$ws.Range("G141").Value = "dnasr281@gmail.com, System"
$ws.Range("G142").Value = "dnasr281@gmail.com, admin@admin.com"
$ws.Range("G145").Value = "dnasr281@gmail.com, System"
$ws.Range("G147").Value = "dnasr281@gmail.com, System"
$ws.Range("G148").Value = "dnasr281@gmail.com, System"
